# Generate Report for Archive
#
# 1) The localization status for the handed-off file moved from
#    "Ready for handoff" to "In Translation" - update the Status cell on
#    every sheet that shows it (the per-locale "zh-cn"/"de-de" sheets keep
#    it in column C, the roll-up "Overview" sheet keeps it once per locale
#    column in row 2).
# 2) The two locale-status columns on the Overview sheet, and the Status
#    column on each locale sheet, got narrower (report archive generation
#    shrank the column to fit the new, shorter status text).

$wb = $excel.ActiveWorkbook

# --- 1) Status text: "Ready for handoff" -> "In Translation" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- 2) Narrow the Status-related columns ---
# Overview: columns E (zh-cn) and F (de-de)
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de: column C (Status)
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
